# Update optimisation data for fasttree: fill in the 7.x rows and add new 8.x rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 56 (7.1) ---
$ws.Range("B56").Value = 44686.05
$ws.Range("C56").Value = -381814.028
$ws.Range("D56").Value = 17263724
$ws.Range("J56").Value = "/usr/bin/time -o 7.1.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 1.6 -close 0.75 -nosupport -nt -gamma global.fa > 7.1.tree"

# --- Row 57 (7.2) ---
$ws.Range("B57").Value = 41110.71
$ws.Range("C57").Value = -381863.763
$ws.Range("D57").Value = 17216428
$ws.Range("J57").Value = "/usr/bin/time -o 7.2.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 1.8 -close 0.75 -nosupport -nt -gamma global.fa > 7.2.tree"

# --- Row 58 (7.3) ---
$ws.Range("B58").Value = 43054.91
$ws.Range("C58").Value = -381819.008
$ws.Range("D58").Value = 17215656
$ws.Range("J58").Value = "/usr/bin/time -o 7.3.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.0 -close 0.75 -nosupport -nt -gamma global.fa > 7.3.tree"

# --- Row 59 (7.4) ---
$ws.Range("B59").Value = 42813.18
$ws.Range("C59").Value = -381604.049
$ws.Range("D59").Value = 17249048
$ws.Range("G59").Value = "right"
$ws.Range("J59").Value = "/usr/bin/time -o 7.4.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -nosupport -nt -gamma global.fa > 7.4.tree"

# --- Row 60 (7.5) ---
$ws.Range("B60").Value = 45508.95
$ws.Range("C60").Value = -381838.781
$ws.Range("D60").Value = 17030732
$ws.Range("J60").Value = "/usr/bin/time -o 7.5.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.5 -close 0.75 -nosupport -nt -gamma global.fa > 7.5.tree"

# --- Row 61 (8.1, new) ---
$ws.Range("A61").Value = 8.1
$ws.Range("B61").Value = 49349.25
$ws.Range("C61").Value = -381604.049
$ws.Range("D61").Value = 17249104

# --- Row 62 (8.2, new) ---
$ws.Range("A62").Value = 8.2
$ws.Range("B62").Value = 59482.7
$ws.Range("C62").Value = -381727.151
$ws.Range("D62").Value = 17283796

# --- Row 63 (8.3, new) ---
$ws.Range("A63").Value = 8.3
$ws.Range("B63").Value = 44243.49
$ws.Range("C63").Value = -381896.94
$ws.Range("D63").Value = 17118112

# --- Row 64 (8.4, new) ---
$ws.Range("A64").Value = 8.4
$ws.Range("B64").Value = 61263.31
$ws.Range("C64").Value = -381815.42
$ws.Range("D64").Value = 17178212

# --- Row 65 (8.5, new) ---
$ws.Range("A65").Value = 8.5
$ws.Range("B65").Value = 54436.22
$ws.Range("C65").Value = -381714.712
$ws.Range("D65").Value = 17648108

# The commandline notes (column J) were typed in this order in the original
# edit: 8.3 first, then 8.1, 8.2, 8.4, 8.5 - reproduce that so the shared
# string table indices line up with the target workbook.
$ws.Range("J63").Value = "/usr/bin/time -o 8.3.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -2nd -nosupport -nt -gamma global.fa > 8.3.tree"
$ws.Range("J61").Value = "/usr/bin/time -o 8.1.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -nt -gamma global.fa > 8.1.tree"
$ws.Range("J62").Value = "/usr/bin/time -o 8.2.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -bionj -nosupport -nt -gamma global.fa > 8.2.tree"
$ws.Range("J64").Value = "/usr/bin/time -o 8.4.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -bionj -2nd -nosupport -nt -gamma global.fa > 8.4.tree"
$ws.Range("J65").Value = "/usr/bin/time -o 8.5.mem.txt -v fasttree -sprlength 20 -refresh 0.8 -topm 2.25 -close 0.75 -fastest -bionj -nosupport -nt -gamma global.fa > 8.5.tree"

# Fill the deltalnL (F) shared formula down to the new rows 56:65
$ws.Range("F56:F65").Formula = '=C56-$C$2'

# Apply the "terminal output" style (Menlo font, as used in rows 3-55) to the
# newly populated data cells by copying formatting from an existing styled cell.
# (Done as separate calls since this runtime does not support multi-area
# union ranges for PasteSpecial.)
$ws.Range("B45").Copy()
$ws.Range("B56:D65").PasteSpecial(-4122)
$ws.Range("B45").Copy()
$ws.Range("J63").PasteSpecial(-4122)
$ws.Range("B45").Copy()
$ws.Range("D66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection/active cell to match where editing finished.
$ws.Range("D66").Select()
